# Debit note template: add "Air", "Telefon" and "Others" as extra values on
# the "Tipe" lookup sheet (rows 5-7, column B) so the W5:W6 dropdown list on
# Sheet1 can offer them too.

$wb = $excel.ActiveWorkbook

$wsTipe = $wb.Worksheets.Item("Tipe")
$wsTipe.Range("B5").Value = "Air"
$wsTipe.Range("B6").Value = "Telefon"
$wsTipe.Range("B7").Value = "Others"

$wsMain = $wb.Worksheets.Item("Sheet1")

# Restore the cursor/selection on each sheet. Select the non-active sheet's
# cell first, then the main sheet's cell last so Sheet1 stays the active
# (tab-selected) sheet, matching the original workbook's layout.
$wsTipe.Range("D9").Select() | Out-Null
$wsMain.Range("V3").Select() | Out-Null
